$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A2:C60")
$sortRange.Sort($ws.Range("A2:A60"), 1)
$sortRange.Select()
